$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new (blank) column before column N, shifting the existing
# "Late" / "heading" / "Outstanding" columns one place to the right.
$ws.Columns("N:N").Insert()
$ws.Columns("N:N").ColumnWidth = 10.14

# Make "Repayment schedule" the active sheet / tab (was "Edit Repayment
# Schedule"), and move the selection to R11.
$ws.Activate()
$ws.Range("R11").Select()
